# Add codes to test what criteria an uncertainty parameter affects.
# The note text referencing "_sys_simulation.py" is corrected to
# "sys_simulation.py" (drop the leading underscore) across the
# DataSummary sheet's G column notes, and the selection is left on the
# last-edited cell to match where the author ended up working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSummary")
$ws.Activate()

$rowsToFix = @(13, 14, 15, 16, 18, 19, 20, 21, 22, 24, 26)
foreach ($r in $rowsToFix) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current -replace "_sys_simulation\.py", "sys_simulation.py"
    }
}

# Reflect the selection recorded for the sheet after the edit.
[void]$ws.Range("G22").Select()
